# Adds the "PC data / PCALLOCATION" data-loading section (including the
# SQL Script block) to the end of the document, right before the final
# (bookmarked) paragraph.

$d = $word.ActiveDocument

# Unicode curly quotes used by a couple of the SQL literals.
$lq = [char]0x2018   # ‘
$rq = [char]0x2019   # ’

# The very last paragraph in the document body holds the _GoBack bookmark.
# All new paragraphs get inserted immediately before it; it keeps its
# bookmark, but gains the final two runs of text plus the "tight" spacing
# that the rest of the new block uses.
$lastIndex = $d.Paragraphs.Count
$anchorRange = $d.Paragraphs.Item($lastIndex).Range

# ---------------------------------------------------------------------
# Plain (non-red, normal-spacing) paragraphs that make up the narrative
# text before the SQL script.
# ---------------------------------------------------------------------
$plainTexts = @(
    "------",
    "The next section of detail data loading that will be shown is how we successfully managed to import the PC data into the PCALLOCATION entity. This caused a large amount of initial confusion due to how the source data was structured in the Excel document. The fact that the numbers for the PC identification numbers changed every time you performed an action in the Excel document meant that there were very real possibilities that there would be duplicates. This of course would pose a major problem especially if they were going to be used as a unique Primary Key. ",
    "Apart from using the Excel function to check for duplicates, another method that we used was first of all enable the Primary Key constraint on the PCID column of the PCALLOCATION entity, and then to continue trying to import them until we come across a combination of numbers that does not have duplicates in them, a crude method to be sure, but effective. ",
    "First of all though, we had to figure a way to import the PC numbers while keeping them linked to their roles. ",
    "The method that we used to tackle this problem was to simply import one column of the PC ids at a time (due to the fact that each column of PC ids had a different role associated with it). Doing this, we were able to ensure that the numbers were unique (by checking to see if there were any violations of the Primary Key at any stage), and also make sure the roles for each PC are accurate by using the following SQL statements at each corresponding stage of the loading.",
    "",
    "SQL Script:"
)

foreach ($t in $plainTexts) {
    $anchorRange.InsertParagraphBefore()
    $newPara = $d.Paragraphs.Item($lastIndex)
    if ($t -ne "") {
        $newPara.Range.Text = $t
    }
    $lastIndex = $d.Paragraphs.Count
    $anchorRange = $d.Paragraphs.Item($lastIndex).Range
}

# ---------------------------------------------------------------------
# The red, tight-line-spacing SQL script block.
# ---------------------------------------------------------------------
$sqlTexts = @(
    "UPDATE PCALLOCATION",
    ("SET PCROLE = " + $lq + "Lectern" + $rq),
    "WHERE PCROLE IS null;",
    "",
    "UPDATE PCALLOCATION",
    "SET PCROLE = 'Office'",
    "WHERE PCROLE IS null;",
    "",
    "UPDATE PCALLOCATION",
    ("SET PCROLE = " + $lq + "Student" + $rq),
    "WHERE PCROLE IS null;"
)

foreach ($t in $sqlTexts) {
    $anchorRange.InsertParagraphBefore()
    $newPara = $d.Paragraphs.Item($lastIndex)
    if ($t -ne "") {
        $newPara.Range.Text = $t
    }
    $newPara.SpaceAfter = 0
    $newPara.LineSpacingRule = 0
    $newPara.Range.Font.Color = 255
    $lastIndex = $d.Paragraphs.Count
    $anchorRange = $d.Paragraphs.Item($lastIndex).Range
}

# ---------------------------------------------------------------------
# Closing blank (tight-spacing, no colour) separator paragraph.
# ---------------------------------------------------------------------
$anchorRange.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($lastIndex)
$newPara.SpaceAfter = 0
$newPara.LineSpacingRule = 0
$lastIndex = $d.Paragraphs.Count
$anchorRange = $d.Paragraphs.Item($lastIndex).Range

# ---------------------------------------------------------------------
# Final paragraph: this is the original bookmark paragraph. It picks up
# the same tight spacing as above, plus the closing sentence, while
# keeping the _GoBack bookmark right where it was (at the very end).
# ---------------------------------------------------------------------
$finalPara = $d.Paragraphs.Item($lastIndex)
$finalPara.SpaceAfter = 0
$finalPara.LineSpacingRule = 0
$anchorRange.InsertBefore("By using each of these SQL statements in turn, we were able to keep the roles of the PCs accurate with the corresponding PC. ")

Write-Output "Done. Paragraph count now: $($d.Paragraphs.Count)"
